$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '65.837.98'
$ws.Cells.Item(2, 5).Value = '  +1.34%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.176.81'
$ws.Cells.Item(3, 5).Value = '  +0.70%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''1.00'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.09%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''594.43'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +3.63%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''152.57'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +1.82%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.999'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.10%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '3.173.13'
$ws.Cells.Item(8, 5).Value = '  +0.67%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.534'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +1.50%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -0.84%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''6.06'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.49%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +2.81%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +0.09%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''38.70'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +4.24%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.695.57'
$ws.Cells.Item(15, 5).Value = '  +0.58%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '65.883.94'
$ws.Cells.Item(16, 5).Value = '  +1.26%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +4.37%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.175.04'
$ws.Cells.Item(18, 5).Value = '  +0.63%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.47%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''507.16'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.24%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''15.31'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +3.02%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''0.733'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +2.17%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''7.98'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +3.43%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''14.94'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -2.42%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''84.61'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.36%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.01%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +3.69%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +2.29%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +5.27%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +12.98%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +3.67%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''27.99'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +1.46%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +2.53%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.02%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -0.82%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''54.72'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -0.05%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''0.0898'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -0.04%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''479.57'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +3.15%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -0.92%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +1.39%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Kaspa'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(41, 4).Value = '''0.121'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +3.31%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'TheGraph'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(42, 4).Value = '''0.298'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +5.50%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'dogwifhat'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(43, 4).Value = '''2.83'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -5.15%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '0.0₃0647'
$ws.Cells.Item(44, 5).Value = '  +10.52%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '2.881.43'
$ws.Cells.Item(45, 5).Value = '  -5.54%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -1.52%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -0.79%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''1.00'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +0.02%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +1.57%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +2.39%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''2.60'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +6.89%  '
